{"js": "// Replace the title date and every two-digit multiplication answer cell.\n// Each \"old\" text is unique within the document, so an exact, case-sensitive\n// search-and-replace for each pair reproduces the diff precisely.\nconst replacements = [\n  [\n    \"2024-08-02 Friday\",\n    \"2024-08-03 Saturday\"\n  ],\n  [\n    \"27\u00d765=1755\",\n    \"69\u00d726=1794\"\n  ],\n  [\n    \"98\u00d711=1078\",\n    \"89\u00d733=2937\"\n  ],\n  [\n    \"37\u00d747=1739\",\n    \"47\u00d727=1269\"\n  ],\n  [\n    \"88\u00d721=1848\",\n    \"38\u00d735=1330\"\n  ],\n  [\n    \"91\u00d727=2457\",\n    \"60\u00d769=4140\"\n  ],\n  [\n    \"40\u00d730=1200\",\n    \"77\u00d719=1463\"\n  ],\n  [\n    \"79\u00d763=4977\",\n    \"33\u00d787=2871\"\n  ],\n  [\n    \"84\u00d797=8148\",\n    \"72\u00d728=2016\"\n  ],\n  [\n    \"57\u00d728=1596\",\n    \"40\u00d771=2840\"\n  ],\n  [\n    \"94\u00d787=8178\",\n    \"64\u00d799=6336\"\n  ],\n  [\n    \"65\u00d718=1170\",\n    \"17\u00d775=1275\"\n  ],\n  [\n    \"20\u00d772=1440\",\n    \"78\u00d744=3432\"\n  ],\n  [\n    \"73\u00d788=6424\",\n    \"64\u00d781=5184\"\n  ],\n  [\n    \"59\u00d793=5487\",\n    \"93\u00d727=2511\"\n  ],\n  [\n    \"99\u00d725=2475\",\n    \"74\u00d780=5920\"\n  ],\n  [\n    \"78\u00d747=3666\",\n    \"84\u00d724=2016\"\n  ],\n  [\n    \"30\u00d798=2940\",\n    \"87\u00d722=1914\"\n  ],\n  [\n    \"41\u00d776=3116\",\n    \"52\u00d771=3692\"\n  ],\n  [\n    \"33\u00d762=2046\",\n    \"34\u00d729=986\"\n  ],\n  [\n    \"32\u00d790=2880\",\n    \"85\u00d790=7650\"\n  ],\n  [\n    \"65\u00d741=2665\",\n    \"94\u00d748=4512\"\n  ],\n  [\n    \"89\u00d786=7654\",\n    \"72\u00d714=1008\"\n  ],\n  [\n    \"71\u00d725=1775\",\n    \"70\u00d751=3570\"\n  ],\n  [\n    \"65\u00d794=6110\",\n    \"65\u00d781=5265\"\n  ],\n  [\n    \"94\u00d785=7990\",\n    \"38\u00d751=1938\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the title date and every two-digit multiplication answer cell.\n# Each \"old\" text is unique within the document, so an exact, case-sensitive\n# Find/Replace for each pair reproduces the diff precisely.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-02 Friday\", \"2024-08-03 Saturday\"),\n    @(\"27\u00d765=1755\", \"69\u00d726=1794\"),\n    @(\"98\u00d711=1078\", \"89\u00d733=2937\"),\n    @(\"37\u00d747=1739\", \"47\u00d727=1269\"),\n    @(\"88\u00d721=1848\", \"38\u00d735=1330\"),\n    @(\"91\u00d727=2457\", \"60\u00d769=4140\"),\n    @(\"40\u00d730=1200\", \"77\u00d719=1463\"),\n    @(\"79\u00d763=4977\", \"33\u00d787=2871\"),\n    @(\"84\u00d797=8148\", \"72\u00d728=2016\"),\n    @(\"57\u00d728=1596\", \"40\u00d771=2840\"),\n    @(\"94\u00d787=8178\", \"64\u00d799=6336\"),\n    @(\"65\u00d718=1170\", \"17\u00d775=1275\"),\n    @(\"20\u00d772=1440\", \"78\u00d744=3432\"),\n    @(\"73\u00d788=6424\", \"64\u00d781=5184\"),\n    @(\"59\u00d793=5487\", \"93\u00d727=2511\"),\n    @(\"99\u00d725=2475\", \"74\u00d780=5920\"),\n    @(\"78\u00d747=3666\", \"84\u00d724=2016\"),\n    @(\"30\u00d798=2940\", \"87\u00d722=1914\"),\n    @(\"41\u00d776=3116\", \"52\u00d771=3692\"),\n    @(\"33\u00d762=2046\", \"34\u00d729=986\"),\n    @(\"32\u00d790=2880\", \"85\u00d790=7650\"),\n    @(\"65\u00d741=2665\", \"94\u00d748=4512\"),\n    @(\"89\u00d786=7654\", \"72\u00d714=1008\"),\n    @(\"71\u00d725=1775\", \"70\u00d751=3570\"),\n    @(\"65\u00d794=6110\", \"65\u00d781=5265\"),\n    @(\"94\u00d785=7990\", \"38\u00d751=1938\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
